# StateFunction.xlsx - "remove Gamelogic project, modify SLG building config"
#
# Sheet1 is a state/function availability matrix: column A holds the
# Func1..Func12 state names (rows 2-13), and columns B..O hold 1/0 flags
# for each EFT_* (event/effect) column. The Gamelogic project column(s)
# are being disabled (flag flipped from 1 to 0) for most Func rows, while
# a handful of cells that are NOT part of the removed project are left
# untouched (e.g. row 2 keeps C2/D2, row 4 keeps D4/F4, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> list of column letters whose flag must become 0.
$changes = @{
    2  = @("E","F","G","H","I","J","K","L","M","N")
    3  = @("C","D","F","G","H","I","J","K","L","M","N")
    4  = @("C","E","G","H","I","J","K","L","M","N")
    5  = @("C","D","F","G","H","I","J","K","L","M","N")
    6  = @("C","D","E","F","G","H","I","J","K","L","M","N")
    7  = @("C","D","E","F","G","H","I","J","K","L","M","N")
    8  = @("C","D","E","F","G","H","I","J","K","L","M","N")
    9  = @("C","D","E","F","G","H","I","J","K","L","M","N")
    10 = @("C","D","E","F","G","H","I","J","K","L","M","N")
    11 = @("C","D","E","F","G","H","I","J","K","L","M","N")
    12 = @("C","D","E","F","G","H","I","J","K","L","M","N")
    13 = @("C","D","E","F","G","H","I","J","K","L","M","N")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 0
    }
}

# The author's last selection in the sheet ended up on F11.
$ws.Range("F11").Select()
